# Fall 2022 refresh: the "updated automatically" date placeholder that
# lives on the slide master and on every slide layout was re-cached by
# PowerPoint from 4/20/22 to 8/9/22 (ppPlaceholderDate, type 16).
# Walk the slide master + every custom (slide) layout, find the date
# placeholder on each one, and push the new cached date into it.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Placeholders.Count; $i++) {
        $ph = $shapes.Placeholders.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = "8/9/22"
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
